$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 45681.978125000001
$ws.Range("B7").Value = 8
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = 188
$ws.Range("E7").Value = 368
$ws.Range("F7").Value = 344
$ws.Range("G7").Value = 384
$ws.Range("H7").Value = 2664
$ws.Range("I7").Value = 384
$ws.Range("J7").Value = 1216
$ws.Range("K7").Value = 119
$ws.Range("L7").Value = 304
$ws.Range("M7").Value = 30
$ws.Range("N7").Value = 2884
$ws.Range("O7").Value = 3651
